# Weekly update: insert two new price records (row 33 = "Primera",
# row 34 = "Segunda") at the top of the Mango data block for the
# Feria Lagunitas de Puerto Montt feed. All existing data rows 33..126
# shift down by two rows (becoming rows 35..128).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 33 (shifts rows 33-126 down to 35-128)
$ws.Rows.Item(33).EntireRow.Insert()
$ws.Rows.Item(34).EntireRow.Insert()

# New row 33 - "Primera" quality entry for the latest reporting date
$ws.Cells.Item(33, 1).Value = 4
$ws.Cells.Item(33, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(33, 3).Value = "Los Lagos"
$ws.Cells.Item(33, 4).Value = 44544
$ws.Cells.Item(33, 5).Value = 10
$ws.Cells.Item(33, 6).Value = "Fruta"
$ws.Cells.Item(33, 7).Value = 100108
$ws.Cells.Item(33, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(33, 9).Value = 100108002
$ws.Cells.Item(33, 10).Value = "Mango"
$ws.Cells.Item(33, 11).Value = "Sin especificar"
$ws.Cells.Item(33, 12).Value = "Primera"
$ws.Cells.Item(33, 13).Value = 200
$ws.Cells.Item(33, 14).Value = 7500
$ws.Cells.Item(33, 15).Value = 8000
$ws.Cells.Item(33, 16).Value = 7750
$ws.Cells.Item(33, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(33, 18).Value = "Perú"
$ws.Cells.Item(33, 19).Value = 1938
$ws.Cells.Item(33, 20).Value = 4

# New row 34 - "Segunda" quality entry for the latest reporting date
$ws.Cells.Item(34, 1).Value = 4
$ws.Cells.Item(34, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(34, 3).Value = "Los Lagos"
$ws.Cells.Item(34, 4).Value = 44544
$ws.Cells.Item(34, 5).Value = 10
$ws.Cells.Item(34, 6).Value = "Fruta"
$ws.Cells.Item(34, 7).Value = 100108
$ws.Cells.Item(34, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(34, 9).Value = 100108002
$ws.Cells.Item(34, 10).Value = "Mango"
$ws.Cells.Item(34, 11).Value = "Sin especificar"
$ws.Cells.Item(34, 12).Value = "Segunda"
$ws.Cells.Item(34, 13).Value = 100
$ws.Cells.Item(34, 14).Value = 5000
$ws.Cells.Item(34, 15).Value = 5000
$ws.Cells.Item(34, 16).Value = 5000
$ws.Cells.Item(34, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(34, 18).Value = "Perú"
$ws.Cells.Item(34, 19).Value = 1250
$ws.Cells.Item(34, 20).Value = 4
